$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, newValue) updates to column F
$updates = @{
    "展览" = @(
        @{Row = 2;  Value = 347}
        @{Row = 3;  Value = 3521}
        @{Row = 5;  Value = 8237}
        @{Row = 7;  Value = 93}
        @{Row = 12; Value = 1199}
        @{Row = 13; Value = 61}
        @{Row = 15; Value = 21}
        @{Row = 16; Value = 587}
        @{Row = 17; Value = 83}
        @{Row = 18; Value = 4574}
        @{Row = 20; Value = 7288}
        @{Row = 22; Value = 55980}
        @{Row = 23; Value = 55980}
        @{Row = 24; Value = 4489}
        @{Row = 26; Value = 1041}
        @{Row = 29; Value = 92}
        @{Row = 32; Value = 3729}
        @{Row = 34; Value = 49}
        @{Row = 37; Value = 1220}
        @{Row = 38; Value = 1194}
        @{Row = 39; Value = 159}
        @{Row = 40; Value = 193}
        @{Row = 41; Value = 1071}
        @{Row = 45; Value = 169}
        @{Row = 46; Value = 11}
    )
    "演出" = @(
        @{Row = 12; Value = 114}
        @{Row = 15; Value = 173}
        @{Row = 16; Value = 7481}
        @{Row = 17; Value = 107}
        @{Row = 28; Value = 33}
        @{Row = 35; Value = 32}
    )
    "本地生活" = @(
        @{Row = 4;  Value = 2293}
        @{Row = 5;  Value = 1553}
        @{Row = 8;  Value = 2339}
        @{Row = 9;  Value = 9337}
        @{Row = 10; Value = 1672}
        @{Row = 11; Value = 162}
        @{Row = 15; Value = 172}
    )
    "全部类型" = @(
        @{Row = 2;  Value = 347}
        @{Row = 3;  Value = 3521}
        @{Row = 5;  Value = 8237}
        @{Row = 6;  Value = 1553}
        @{Row = 8;  Value = 2339}
        @{Row = 9;  Value = 1672}
        @{Row = 10; Value = 162}
        @{Row = 13; Value = 93}
        @{Row = 16; Value = 61}
        @{Row = 17; Value = 21}
        @{Row = 18; Value = 587}
        @{Row = 19; Value = 83}
        @{Row = 20; Value = 7288}
        @{Row = 21; Value = 55980}
        @{Row = 24; Value = 4489}
        @{Row = 25; Value = 1041}
        @{Row = 27; Value = 92}
        @{Row = 29; Value = 114}
        @{Row = 30; Value = 3729}
        @{Row = 32; Value = 49}
        @{Row = 35; Value = 1220}
        @{Row = 36; Value = 107}
        @{Row = 37; Value = 159}
        @{Row = 38; Value = 193}
        @{Row = 39; Value = 1071}
        @{Row = 42; Value = 169}
        @{Row = 47; Value = 32}
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($entry in $updates[$sheetName]) {
        $ws.Cells.Item($entry.Row, 6).Value = $entry.Value
    }
}
